$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 2.86
    "J2" = 3.25
    "N2" = 3.5
    "T2" = 1.73
    "Z2" = 17.5
    "AO2" = 24
    "F3" = 1.84
    "G3" = 1.85
    "S3" = 2.56
    "W3" = 2.18
    "Y3" = 23
    "AA3" = 95
    "AE3" = 48
    "AI3" = 48
    "AJ3" = 20
    "AO3" = 38
    "F4" = 9.6
    "G4" = 12
    "H4" = 1.33
    "J4" = 5.7
    "R4" = 1.54
    "T4" = 2
    "W4" = 1.09
    "AE4" = 17.5
    "AO4" = 6
    "N5" = 1.01
    "P5" = 2.32
    "T5" = 1.59
    "H6" = 9.199999999999999
    "P6" = 2.94
    "S6" = 2.22
    "T6" = 1.77
    "U6" = 2.22
    "Y6" = 42
    "Z6" = 90
    "AA6" = 290
    "G7" = 3.7
    "V7" = 1.87
    "W7" = 1.37
    "AN7" = 28
    "L8" = 1.4
    "N8" = 3.8
    "P8" = 1.95
    "R8" = 1.37
    "S8" = 3.5
    "T8" = 2
    "U8" = 1.96
    "X8" = 13.5
    "AB8" = 8.199999999999999
    "AE8" = 90
    "AI8" = 90
    "AN8" = 10.5
    "G12" = 2.16
    "H12" = 3.85
    "J12" = 3.45
    "K12" = 3.9
    "N12" = 1.86
    "P12" = 1.86
    "Q12" = 1.96
    "S12" = 2.98
    "T12" = 1.68
    "W12" = 1.86
    "X12" = 19.5
    "Z12" = 42
    "AB12" = 12.5
    "AC12" = 11
    "AD12" = 23
    "AF12" = 17.5
    "G13" = 2.76
    "H13" = 2.62
    "O13" = 1.2
    "Q13" = 1.51
    "R13" = 1.57
    "S13" = 2.54
    "F14" = 1.65
    "G14" = 1.76
    "H14" = 4.9
    "I14" = 5.7
    "M14" = 1.04
    "N14" = 4.8
    "O14" = 1.22
    "R14" = 1.53
    "S14" = 2.56
    "T14" = 1.69
    "U14" = 2.18
    "V14" = 1.22
    "X14" = 24
    "Y14" = 23
    "Z14" = 44
    "AA14" = 140
    "AB14" = 11.5
    "AC14" = 11
    "AD14" = 21
    "AE14" = 60
    "AF14" = 12.5
    "AG14" = 11
    "AH14" = 19.5
    "AI14" = 60
    "AJ14" = 18.5
    "AK14" = 17.5
    "AL14" = 32
    "AM14" = 100
    "AN14" = 8.199999999999999
    "AO14" = 65
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"